{"js": "// Apply the three text replacements described by the diff.\n// Each entry is [oldText, newText]; we locate the paragraph containing the\n// old text (via Body.search, an exact, case-sensitive match) and replace\n// the whole run's text in place so surrounding paragraph formatting is kept.\nconst replacements = [\n  {\n    oldText:\n      \"\u8981\u6c42\u672c\u7cfb\u7edf\u5728\u5b8c\u6210\u5404\u9879\u529f\u80fd\u7684\u540c\u65f6\uff0c\u8981\u6c42\u7cfb\u7edf\u5904\u7406\u8fc5\u901f\uff0c\u5904\u7406\u4e8b\u52a1\u9700\u8981\u957f\u65f6\u95f4\u65f6\uff0c\u63d0\u793a\u7528\u6237\u7b49\u5f85\u4e14\u7b49\u5f85\u65f6\u95f4\u5728\u7528\u6237\u53ef\u63a5\u6536\u7684\u8303\u56f4\u4e4b\u5185\u3002\",\n    newText:\n      \"\u8981\u6c42\u672c\u7cfb\u7edf\u5728\u5c40\u57df\u7f51\u8303\u56f4\u5185\u4f7f\u7528\uff0c\u4e14\u540c\u65f6\u4f7f\u7528\u4eba\u6570\u572810\u4eba\u4ee5\u5185\u65f6\uff0c\u4e0d\u4ec5\u53ef\u4ee5\u5b8c\u6210\u5168\u90e8\u7684\u529f\u80fd\uff0c\u8fd8\u80fd\u5904\u7406\u8fc5\u901f\u3002\u4e00\u822c\u5728\u7528\u6237\u8fdb\u884c\u76f8\u5e94\u64cd\u4f5c\u76843\u79d2\u5185\uff0c\u5c31\u80fd\u8f93\u51fa\u6240\u9700\u7684\u5904\u7406\u7ed3\u679c\u3002\u5904\u7406\u4e8b\u52a1\u9700\u8981\u957f\u65f6\u95f4\u65f6\uff0c\u63d0\u793a\u7528\u6237\u7b49\u5f85\u4e14\u7b49\u5f85\u65f6\u95f4\u572810\u79d2\u4e4b\u5185\u3002\",\n  },\n  {\n    oldText:\n      \"\u8fd9\u9879\u529f\u80fd\u4e3b\u8981\u662f\u4e3a\u4e86\u7528\u6237\u8fdb\u884c\u67d0\u4e9b\u64cd\u4f5c\u540e\uff0c\u6ca1\u6709\u4fdd\u5b58\u5230\u6587\u4ef6\u524d\uff0c\u53ef\u4ee5\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u8fd9\u4e9b\u6570\u636e\u4fe1\u606f\uff0c\u4f7f\u7528\u6237\u7684\u64cd\u4f5c\u4e0d\u8d77\u4f5c\u7528\u3002\",\n    newText:\n      \"\u8fd9\u9879\u529f\u80fd\u4e3b\u8981\u662f\u4e3a\u4e86\u7528\u6237\u8fdb\u884c\u67d0\u4e9b\u64cd\u4f5c\u540e\uff0c\u6ca1\u6709\u4fdd\u5b58\u5230\u6570\u636e\u5e93\u524d\uff0c\u53ef\u4ee5\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u8fd9\u4e9b\u6570\u636e\u4fe1\u606f\uff0c\u4f7f\u7528\u6237\u7684\u64cd\u4f5c\u4e0d\u8d77\u4f5c\u7528\u3002\",\n  },\n  {\n    oldText:\n      \"\u5f53\u7528\u6237\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e\uff0c\u53ef\u9000\u51fa\u8be5\u7cfb\u7edf\u3002\u5982\u679c\u7528\u6237\u8fdb\u884c\u4e86\u5f71\u54cd\u5de5\u7a0b\u5e08\u8d44\u6599\u4fe1\u606f\u7684\u64cd\u4f5c\uff0c\u63d0\u793a\u7528\u6237\u662f\u5426\u8fdb\u884c\u4fdd\u5b58\u3002\",\n    newText:\n      \"\u5f53\u7528\u6237\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e\uff0c\u53ef\u8fdb\u884c\u6ce8\u9500\u64cd\u4f5c\u3002\u5982\u679c\u7528\u6237\u8fdb\u884c\u4e86\u5f71\u54cd\u5de5\u7a0b\u5e08\u8d44\u6599\u4fe1\u606f\u7684\u64cd\u4f5c\uff0c\u63d0\u793a\u7528\u6237\u662f\u5426\u5728\u6ce8\u9500\u6216\u79bb\u5f00\u7f51\u9875\u524d\u8fdb\u884c\u4fdd\u5b58\u3002\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target text: \" + oldText);\n  }\n\n  // Replace the exact matched range's text, preserving its formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the three text replacements described by the diff using\n# Word's Find/Replace (wdReplaceAll) against the whole document content.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"\u8981\u6c42\u672c\u7cfb\u7edf\u5728\u5b8c\u6210\u5404\u9879\u529f\u80fd\u7684\u540c\u65f6\uff0c\u8981\u6c42\u7cfb\u7edf\u5904\u7406\u8fc5\u901f\uff0c\u5904\u7406\u4e8b\u52a1\u9700\u8981\u957f\u65f6\u95f4\u65f6\uff0c\u63d0\u793a\u7528\u6237\u7b49\u5f85\u4e14\u7b49\u5f85\u65f6\u95f4\u5728\u7528\u6237\u53ef\u63a5\u6536\u7684\u8303\u56f4\u4e4b\u5185\u3002\"\n        New = \"\u8981\u6c42\u672c\u7cfb\u7edf\u5728\u5c40\u57df\u7f51\u8303\u56f4\u5185\u4f7f\u7528\uff0c\u4e14\u540c\u65f6\u4f7f\u7528\u4eba\u6570\u572810\u4eba\u4ee5\u5185\u65f6\uff0c\u4e0d\u4ec5\u53ef\u4ee5\u5b8c\u6210\u5168\u90e8\u7684\u529f\u80fd\uff0c\u8fd8\u80fd\u5904\u7406\u8fc5\u901f\u3002\u4e00\u822c\u5728\u7528\u6237\u8fdb\u884c\u76f8\u5e94\u64cd\u4f5c\u76843\u79d2\u5185\uff0c\u5c31\u80fd\u8f93\u51fa\u6240\u9700\u7684\u5904\u7406\u7ed3\u679c\u3002\u5904\u7406\u4e8b\u52a1\u9700\u8981\u957f\u65f6\u95f4\u65f6\uff0c\u63d0\u793a\u7528\u6237\u7b49\u5f85\u4e14\u7b49\u5f85\u65f6\u95f4\u572810\u79d2\u4e4b\u5185\u3002\"\n    },\n    @{\n        Old = \"\u8fd9\u9879\u529f\u80fd\u4e3b\u8981\u662f\u4e3a\u4e86\u7528\u6237\u8fdb\u884c\u67d0\u4e9b\u64cd\u4f5c\u540e\uff0c\u6ca1\u6709\u4fdd\u5b58\u5230\u6587\u4ef6\u524d\uff0c\u53ef\u4ee5\u4ece\u6587\u4ef6\u91cd\u65b0\u5f97\u5230\u8fd9\u4e9b\u6570\u636e\u4fe1\u606f\uff0c\u4f7f\u7528\u6237\u7684\u64cd\u4f5c\u4e0d\u8d77\u4f5c\u7528\u3002\"\n        New = \"\u8fd9\u9879\u529f\u80fd\u4e3b\u8981\u662f\u4e3a\u4e86\u7528\u6237\u8fdb\u884c\u67d0\u4e9b\u64cd\u4f5c\u540e\uff0c\u6ca1\u6709\u4fdd\u5b58\u5230\u6570\u636e\u5e93\u524d\uff0c\u53ef\u4ee5\u4ece\u6570\u636e\u5e93\u91cd\u65b0\u5f97\u5230\u8fd9\u4e9b\u6570\u636e\u4fe1\u606f\uff0c\u4f7f\u7528\u6237\u7684\u64cd\u4f5c\u4e0d\u8d77\u4f5c\u7528\u3002\"\n    },\n    @{\n        Old = \"\u5f53\u7528\u6237\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e\uff0c\u53ef\u9000\u51fa\u8be5\u7cfb\u7edf\u3002\u5982\u679c\u7528\u6237\u8fdb\u884c\u4e86\u5f71\u54cd\u5de5\u7a0b\u5e08\u8d44\u6599\u4fe1\u606f\u7684\u64cd\u4f5c\uff0c\u63d0\u793a\u7528\u6237\u662f\u5426\u8fdb\u884c\u4fdd\u5b58\u3002\"\n        New = \"\u5f53\u7528\u6237\u4e0d\u518d\u4f7f\u7528\u8be5\u7cfb\u7edf\u540e\uff0c\u53ef\u8fdb\u884c\u6ce8\u9500\u64cd\u4f5c\u3002\u5982\u679c\u7528\u6237\u8fdb\u884c\u4e86\u5f71\u54cd\u5de5\u7a0b\u5e08\u8d44\u6599\u4fe1\u606f\u7684\u64cd\u4f5c\uff0c\u63d0\u793a\u7528\u6237\u662f\u5426\u5728\u6ce8\u9500\u6216\u79bb\u5f00\u7f51\u9875\u524d\u8fdb\u884c\u4fdd\u5b58\u3002\"\n    }\n)\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.Text = $rep.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$rep.New, 2) | Out-Null\n}\n"}
